# "Fruta / hortaliza, semanal" — weekly refresh of the daily price sheet.
#
# A new daily price record was appended upstream; in this sheet's layout
# (most-recent-first within the "Región de Ñuble" / Primera block) that
# shows up as a brand-new row inserted at row 62, pushing every existing
# record from row 62 down through row 184 to rows 63 through 185
# (dimension grows from A1:R184 to A1:R185).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 62; Excel shifts rows 62:184 down to 63:185 and
# carries the column D date-number formatting onto the new blank row.
$ws.Rows(62).Insert()

# Populate the newly inserted row with the new daily record.
$ws.Range("A62").Value = 11
$ws.Range("B62").Value = "Vega Monumental Concepción"
$ws.Range("C62").Value = "Bíobío"
$ws.Range("D62").Value = 44533
$ws.Range("E62").Value = 8
$ws.Range("F62").Value = 100112009
$ws.Range("G62").Value = "Acelga"
$ws.Range("H62").Value = "Sin especificar"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 350
$ws.Range("K62").Value = 600
$ws.Range("L62").Value = 650
$ws.Range("M62").Value = 621
$ws.Range("N62").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O62").Value = "Región de Ñuble"
$ws.Range("P62").Value = 621
$ws.Range("Q62").Value = 1
$ws.Range("R62").Value = "Hortaliza"
